$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.630.58'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.434.67'
$ws.Range('E3').Value = '  -3.06%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '591.49'
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.88'
$ws.Range('E6').Value = '  -5.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.434.25'
$ws.Range('E7').Value = '  -2.97%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.30'
$ws.Range('E10').Value = '  -6.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.122'
$ws.Range('E11').Value = '  -8.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -7.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.023.30'
$ws.Range('E13').Value = '  -2.92%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000181'
$ws.Range('E14').Value = '  -10.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.48'
$ws.Range('E15').Value = '  -8.84%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.455.71'
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.590.61'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.74'
$ws.Range('E19').Value = '  -11.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.92'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.69'
$ws.Range('E21').Value = '  -6.46%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '393.85'
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.554'
$ws.Range('E23').Value = '  -7.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '73.45'
$ws.Range('E24').Value = '  -5.54%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.587.25'
$ws.Range('E26').Value = '  -2.59%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000107'
$ws.Range('E27').Value = '  -8.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.28'
$ws.Range('E29').Value = '  -9.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.19'
$ws.Range('E30').Value = '  -8.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.24'
$ws.Range('E31').Value = '  -9.34%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.444.42'
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.145'
$ws.Range('E34').Value = '  -7.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '23.02'
$ws.Range('E35').Value = '  -5.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '173.31'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.91'
$ws.Range('E37').Value = '  -8.76%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.17'
$ws.Range('E38').Value = '  -9.60%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.47'
$ws.Range('E39').Value = '  -7.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.80'
$ws.Range('E40').Value = '  -9.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0763'
$ws.Range('E41').Value = '  -7.62%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.825'
$ws.Range('E42').Value = '  -4.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '43.87'
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.44'
$ws.Range('E45').Value = '  -13.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.62'
$ws.Range('E46').Value = '  -9.86%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.98'
$ws.Range('E48').Value = '  +1.24%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.58'
$ws.Range('E49').Value = '  -7.63%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.12'
$ws.Range('E50').Value = '  -13.30%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.213.22'
$ws.Range('E51').Value = '  -6.50%  '
